$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated vm_pu values for case with 380 kV (Case_1_72)
$updates = @{
    "B2" = 1.02
    "C2" = 1.018551740108789
    "D2" = 1.032085247126693
    "E2" = 1.019260231493731
    "F2" = 1.029400682875663
    "I2" = 1.031357216112601
    "J2" = 1.023759529937032
    "K2" = 1.03489149598025
    "L2" = 1.022104006589257
    "M2" = 1.032214702682984
    "N2" = 1.025213385896405
    "B3" = 1.02
    "C3" = 1.02259471592138
    "D3" = 1.032658129095715
    "E3" = 1.022825377187387
    "F3" = 1.03286524420154
    "I3" = 1.03157668253968
    "J3" = 1.027422724284277
    "K3" = 1.035273708230464
    "L3" = 1.025467570921765
    "M3" = 1.035480268625895
    "N3" = 1.028881782399794
    "B4" = 1.02
    "C4" = 1.025180922297191
    "D4" = 1.033026052783916
    "E4" = 1.025104623030046
    "F4" = 1.035079352626473
    "I4" = 1.031714758696313
    "J4" = 1.029764225318075
    "K4" = 1.035517561627465
    "L4" = 1.027616456917203
    "M4" = 1.037565647778943
    "N4" = 1.031226608633634
    "B5" = 1.02
    "C5" = 1.026261261467706
    "D5" = 1.033180078000812
    "E5" = 1.026056420764929
    "F5" = 1.036003739637638
    "I5" = 1.031771877700203
    "J5" = 1.030741918991965
    "K5" = 1.035619261703682
    "L5" = 1.028513462220322
    "M5" = 1.038435918074989
    "N5" = 1.032205690744684
    "B6" = 1.02
    "C6" = 1.02644225750902
    "D6" = 1.033205901772067
    "E6" = 1.026215862923247
    "F6" = 1.03615857742619
    "I6" = 1.031781414182662
    "J6" = 1.030905693239783
    "K6" = 1.035636290126165
    "L6" = 1.02866370490974
    "M6" = 1.038581669289785
    "N6" = 1.032369697570719
    "B7" = 1.02
    "C7" = 1.02519538463801
    "D7" = 1.033028113412901
    "E7" = 1.025117365847573
    "F7" = 1.035091729305032
    "I7" = 1.031715525552902
    "J7" = 1.029777315232447
    "K7" = 1.035518923737348
    "L7" = 1.027628467562771
    "M7" = 1.03757730135234
    "N7" = 1.031239717137186
    "B8" = 1.02
    "C8" = 1.019924419897437
    "D8" = 1.032279436699466
    "E8" = 1.020470947661268
    "F8" = 1.030577408168158
    "I8" = 1.031432208480933
    "J8" = 1.025003633461831
    "K8" = 1.035021391697045
    "L8" = 1.023246576884426
    "M8" = 1.033324165172427
    "N8" = 1.026459256190907
    "B9" = 1.02
    "C9" = 1.010395942710241
    "D9" = 1.030938349212857
    "E9" = 1.012061390676767
    "F9" = 1.022400890610748
    "I9" = 1.030902214528256
    "J9" = 1.016360320787704
    "K9" = 1.034117527640776
    "L9" = 1.015304089129871
    "M9" = 1.025608478093193
    "N9" = 1.01780366902138
    "B10" = 1.02
    "C10" = 1.003865136955737
    "D10" = 1.030028715266497
    "E10" = 1.006290831352109
    "F10" = 1.016786847685157
    "I10" = 1.030527289269304
    "J10" = 1.010426958710374
    "K10" = 1.033495785085654
    "L10" = 1.009846047399647
    "M10" = 1.020302527007828
    "N10" = 1.011861880889335
    "B11" = 1.02
    "C11" = 1.000990925782749
    "D11" = 1.029630924756604
    "E11" = 1.003749662181936
    "F11" = 1.014313956249759
    "I11" = 1.030359602994066
    "J11" = 1.007813505766621
    "K11" = 1.033221801245122
    "L11" = 1.007440590443764
    "M11" = 1.017963360363739
    "N11" = 1.009244716542635
    "B12" = 1.02
    "C12" = 0.9999159993599543
    "D12" = 1.02948256029145
    "E12" = 1.002799061748618
    "F12" = 1.013388813763087
    "I12" = 1.030296494775811
    "J12" = 1.006835774857923
    "K12" = 1.033119295937925
    "L12" = 1.006540468782242
    "M12" = 1.017087947040723
    "N12" = 1.008265597143899
    "B13" = 1.02
    "C13" = 1.000146911721357
    "D13" = 1.029514412781472
    "E13" = 1.003003277056064
    "F13" = 1.01358756354116
    "I13" = 1.030310069229426
    "J13" = 1.007045822667254
    "K13" = 1.033141317283667
    "L13" = 1.006733852913699
    "M13" = 1.017276026963831
    "N13" = 1.008475943245212
    "B14" = 1.02
    "C14" = 1.000902223402188
    "D14" = 1.029618673409837
    "E14" = 1.003671223747875
    "F14" = 1.014237620084408
    "I14" = 1.030354403335992
    "J14" = 1.007732830525806
    "K14" = 1.033213343245516
    "L14" = 1.007366323125655
    "M14" = 1.017891133588194
    "N14" = 1.009163926733724
    "B15" = 1.02
    "C15" = 1.001366615417362
    "D15" = 1.029682830743046
    "E15" = 1.004081870639534
    "F15" = 1.014637257696433
    "I15" = 1.030381609506377
    "J15" = 1.008155184010005
    "K15" = 1.033257622752327
    "L15" = 1.007755121284556
    "M15" = 1.018269245367351
    "N15" = 1.009586880008316
    "B16" = 1.02
    "C16" = 1.004054882794083
    "D16" = 1.030055031231574
    "E16" = 1.006458558854757
    "F16" = 1.016950056200913
    "I16" = 1.030538303870927
    "J16" = 1.010599444721636
    "K16" = 1.03351386644436
    "L16" = 1.01000477724179
    "M16" = 1.020456868204918
    "N16" = 1.012034611850517
    "B17" = 1.02
    "C17" = 1.005728512389167
    "D17" = 1.030287441719541
    "E17" = 1.007937802134501
    "F17" = 1.018389372636284
    "I17" = 1.030635150584599
    "J17" = 1.012120586983406
    "K17" = 1.033673312048767
    "L17" = 1.01140444830275
    "M17" = 1.021817757501127
    "N17" = 1.013557914308777
    "B18" = 1.02
    "C18" = 1.00670026705781
    "D18" = 1.030422626646426
    "E18" = 1.008796543182847
    "F18" = 1.019224872315401
    "I18" = 1.030691125666062
    "J18" = 1.013003595237401
    "K18" = 1.033765855574853
    "L18" = 1.012216813112347
    "M18" = 1.022607542043591
    "N18" = 1.014442176535774
    "B19" = 1.02
    "C19" = 1.00703086683481
    "D19" = 1.030468657980575
    "E19" = 1.009088669553481
    "F19" = 1.019509081389879
    "I19" = 1.030710125123737
    "J19" = 1.013303967326501
    "K19" = 1.033797333351537
    "L19" = 1.012493132327152
    "M19" = 1.022876168247953
    "N19" = 1.014742975187694
    "B20" = 1.02
    "C20" = 1.005549410073717
    "D20" = 1.030262545321156
    "E20" = 1.007779517213012
    "F20" = 1.018235366450709
    "I20" = 1.030624813149593
    "J20" = 1.011957824550456
    "K20" = 1.033656252597007
    "L20" = 1.011254696836417
    "M20" = 1.021672162623242
    "N20" = 1.013394920734502
    "B21" = 1.02
    "C21" = 1.000680007854547
    "D21" = 1.029587988181172
    "E21" = 1.003474717575631
    "F21" = 1.014046379278065
    "I21" = 1.030341370902386
    "J21" = 1.007530719175537
    "K21" = 1.033192153861156
    "L21" = 1.007180261894089
    "M21" = 1.017710182983904
    "N21" = 1.008961528362155
    "B22" = 1.02
    "C22" = 0.9975759015947813
    "D22" = 1.029160342928275
    "E22" = 1.000729204590557
    "F22" = 1.011374251878006
    "I22" = 1.030158392815347
    "J22" = 1.004706675024057
    "K22" = 1.032896090672458
    "L22" = 1.004579994409794
    "M22" = 1.015181124007475
    "N22" = 1.006133473744057
    "B23" = 1.02
    "C23" = 0.9992256009589168
    "D23" = 1.02938738662091
    "E23" = 1.002188450934511
    "F23" = 1.012794533941582
    "I23" = 1.030255851605459
    "J23" = 1.006207711251992
    "K23" = 1.033053450541328
    "L23" = 1.005962201196579
    "M23" = 1.01652552716082
    "N23" = 1.007636641615608
    "B24" = 1.02
    "C24" = 1.005630352412625
    "D24" = 1.030273796093277
    "E24" = 1.007851051931349
    "F24" = 1.018304967647868
    "I24" = 1.030629485778444
    "J24" = 1.012031382993403
    "K24" = 1.033663962444249
    "L24" = 1.011322375540384
    "M24" = 1.021737963020185
    "N24" = 1.013468583638875
    "B25" = 1.02
    "C25" = 1.012889457579419
    "D25" = 1.03128772873276
    "E25" = 1.014263264841156
    "F25" = 1.024542389085463
    "I25" = 1.031042966980155
    "J25" = 1.018623810091621
    "K25" = 1.034354506654266
    "L25" = 1.017385052638039
    "M25" = 1.027630723281262
    "N25" = 1.020070372739736
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells in vm_pu.xlsx (Case_1_72, 380 kV)"
